$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: move existing data from columns A:E to B:F (direct literal writes to avoid
#     precision drift from a generic column-insert/shift operation) ---
$ws.Range("B2").Value = -21.2
$ws.Range("C2").Value = 3.8
$ws.Range("D2").Value = -10.8
$ws.Range("E2").Value = -8.220000000000001
$ws.Range("F2").ClearContents()
$ws.Range("B3").Value = -21.67
$ws.Range("C3").Value = 4.54
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = -6.85
$ws.Range("F3").Value = 13.95
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = -7.43
$ws.Range("F4").Value = 10.77
$ws.Range("B5").Value = -23.4
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = -15.14
$ws.Range("E5").Value = -6.27
$ws.Range("F5").Value = 9.970000000000001
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = 7.87
$ws.Range("D6").Value = -12.89
$ws.Range("E6").Value = -8.59
$ws.Range("F6").Value = 14.65
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = 6.56
$ws.Range("D7").Value = -10.59
$ws.Range("E7").Value = -9.31
$ws.Range("F7").Value = 14.77
$ws.Range("B8").Value = -21.48
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = -11.02
$ws.Range("E8").Value = -7.09
$ws.Range("F8").Value = 12.94
$ws.Range("B9").Value = -19.24
$ws.Range("C9").Value = 10.18
$ws.Range("D9").ClearContents()
$ws.Range("E9").Value = -9.9
$ws.Range("F9").Value = 12.31
$ws.Range("B10").Value = -18.76
$ws.Range("C10").Value = 9.35
$ws.Range("D10").Value = -10.58
$ws.Range("E10").Value = -8.93
$ws.Range("F10").Value = 16.23
$ws.Range("B11").Value = -22.23
$ws.Range("C11").Value = 6.01
$ws.Range("D11").ClearContents()
$ws.Range("E11").Value = -8.01
$ws.Range("F11").Value = 13.91
$ws.Range("B12").Value = -24
$ws.Range("C12").Value = 5.38
$ws.Range("D12").Value = -15.79
$ws.Range("E12").Value = -10.06
$ws.Range("F12").Value = 13.01
$ws.Range("B13").Value = -21.95
$ws.Range("C13").Value = 4.42
$ws.Range("D13").Value = -13.63
$ws.Range("E13").Value = -8.92
$ws.Range("F13").Value = 12.56
$ws.Range("B14").Value = -20.07
$ws.Range("C14").Value = 8.33
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = -7.97
$ws.Range("F14").Value = 12.97
$ws.Range("B15").Value = -21.67
$ws.Range("C15").Value = 5.53
$ws.Range("D15").Value = -10.67
$ws.Range("E15").Value = -6.43
$ws.Range("F15").Value = 14.11
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("D16").Value = -9.779999999999999
$ws.Range("E16").Value = -4.8
$ws.Range("F16").Value = 13.16
$ws.Range("B17").Value = -20.43
$ws.Range("C17").Value = 11.21
$ws.Range("D17").Value = -12.35
$ws.Range("E17").Value = -6.63
$ws.Range("F17").Value = 8.15
$ws.Range("B18").Value = -23.05
$ws.Range("C18").Value = 5.81
$ws.Range("D18").ClearContents()
$ws.Range("E18").Value = -8.869999999999999
$ws.Range("F18").Value = 12.54
$ws.Range("B19").Value = -22.15
$ws.Range("C19").Value = 5.44
$ws.Range("D19").Value = -14.32
$ws.Range("E19").Value = -8.470000000000001
$ws.Range("F19").Value = 13.12
$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = 5.51
$ws.Range("D20").Value = -14.65
$ws.Range("E20").Value = -8.140000000000001
$ws.Range("F20").ClearContents()
$ws.Range("B21").Value = -21.82
$ws.Range("C21").Value = 7.49
$ws.Range("D21").Value = -12.06
$ws.Range("E21").Value = -8.890000000000001
$ws.Range("F21").Value = 13.76
$ws.Range("B22").Value = -20.95
$ws.Range("C22").ClearContents()
$ws.Range("D22").Value = -10.65
$ws.Range("E22").Value = -5.14
$ws.Range("F22").Value = 12.2
$ws.Range("B23").Value = -20.86
$ws.Range("C23").Value = 5.8
$ws.Range("D23").Value = -10.59
$ws.Range("E23").Value = -8.529999999999999
$ws.Range("F23").Value = 10.67
$ws.Range("B24").Value = -20.62
$ws.Range("C24").Value = 9.130000000000001
$ws.Range("D24").Value = -12.7
$ws.Range("E24").Value = -9.24
$ws.Range("F24").Value = 12.16
$ws.Range("B25").Value = -21.66
$ws.Range("C25").Value = 5.09
$ws.Range("D25").ClearContents()
$ws.Range("E25").Value = -7.71
$ws.Range("F25").Value = 15.03

# --- Step 2: write the new header row B1:F1 (shifted from old A1:E1) ---
$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"
$ws.Range("E1").Value = "D"
$ws.Range("F1").Value = "F"

# --- Step 3: add the new "ID" column in column A ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "Hb 2"
$ws.Range("A3").Value = "Hb 3"
$ws.Range("A4").Value = "S 24"
$ws.Range("A5").Value = "S 28"
$ws.Range("A6").Value = "Hb 107"
$ws.Range("A7").Value = "Hb 66"
$ws.Range("A8").Value = "Hb 69"
$ws.Range("A9").Value = "Hb 95"
$ws.Range("A10").Value = "Hb 99"
$ws.Range("A11").Value = "Hb 92"
$ws.Range("A12").Value = "Hb 40"
$ws.Range("A13").Value = "Hb 41"
$ws.Range("A14").Value = "S 11"
$ws.Range("A15").Value = "Hb 57"
$ws.Range("A16").Value = "S 21"
$ws.Range("A17").Value = "S 22"
$ws.Range("A18").Value = "S 3"
$ws.Range("A19").Value = "S 4"
$ws.Range("A20").Value = "S 5"
$ws.Range("A21").Value = "Hb 74"
$ws.Range("A22").Value = "Hb 79"
$ws.Range("A23").Value = "Hb 32"
$ws.Range("A24").Value = "S 15"
$ws.Range("A25").Value = "S 16"
